# Fruta / hortaliza, semanal
# Insert a new record row at row 134 (pushing the existing rows 134-137 down
# to 135-138) in the "Feria Lagunitas de Puerto Montt - Pomelo" sheet, and
# populate the newly inserted row with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 134:137 down to 135:138, inserting a blank row 134.
$ws.Rows.Item(134).Insert()

# Fill in the new row 134 with the new data.
$ws.Range("A134").Value2 = 4
$ws.Range("B134").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value2 = "Los Lagos"
$ws.Range("D134").Value2 = 44448
$ws.Range("E134").Value2 = 10
$ws.Range("F134").Value2 = "Fruta"
$ws.Range("G134").Value2 = 100102
$ws.Range("H134").Value2 = "Cítricos"
$ws.Range("I134").Value2 = 100102006
$ws.Range("J134").Value2 = "Pomelo"
$ws.Range("K134").Value2 = "Start Ruby"
$ws.Range("L134").Value2 = "Primera"
$ws.Range("M134").Value2 = 60
$ws.Range("N134").Value2 = 12000
$ws.Range("O134").Value2 = 12000
$ws.Range("P134").Value2 = 12000
$ws.Range("Q134").Value2 = "`$/caja 14 kilos empedrada"
$ws.Range("R134").Value2 = "Región de O'Higgins"
$ws.Range("S134").Value2 = 857
$ws.Range("T134").Value2 = 14
